$d = $word.ActiveDocument

# Find the paragraph index of the "{ownerAddress}" line (right-aligned block
# sitting right above the "Objet : ..." paragraph).
$ownerAddressIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*{ownerAddress}*") {
        $ownerAddressIdx = $i
        break
    }
}

# Insert a new right-aligned paragraph right after it and fill it with
# "{ownerEmail}" (inherits the {ownerAddress} paragraph's pPr: spacing
# after=0, right-aligned, no Arial override).
$p = $d.Paragraphs.Item($ownerAddressIdx)
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertParagraphAfter()
$emailIdx = $ownerAddressIdx + 1
$d.Paragraphs.Item($emailIdx).Range.Text = "{ownerEmail}"

# Insert another right-aligned paragraph after that one for "{ownerPhoneNumber}".
$p = $d.Paragraphs.Item($emailIdx)
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertParagraphAfter()
$phoneIdx = $emailIdx + 1
$d.Paragraphs.Item($phoneIdx).Range.Text = "{ownerPhoneNumber}"

# Insert one more (empty) right-aligned paragraph after that.
$p = $d.Paragraphs.Item($phoneIdx)
$r = $p.Range.Duplicate
$r.Collapse(0)
$r.InsertParagraphAfter()

# Merge the "Objet : " / "Attestation de remise des clefs" runs into a
# single run (same visible text, but Word's Find/Replace normalises the
# adjacent identically-formatted runs into one).
$d.Content.Find.Execute("Objet : Attestation de remise des clefs", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Objet : Attestation de remise des clefs", 2)

Write-Output "done"
